# Cover Letter update:
#   1. Bump the letter date from March 5, 2021 to March 30, 2021.
#   2. Reword "...within the cutting-edge technological environment..."
#      to "...within a cutting-edge technological environment...".
#
# For (2) the source docx stores the whole sentence as a single run; the
# target keeps that exact text run-for-run except that the word "a" that
# replaces "the" must land in its own <w:r> (with identical rPr) while the
# "Coinbase" / ". " / "As a data analyst..." / "Thank you..." runs that
# follow stay exactly as separate runs too. A plain Find/Replace across the
# paragraph causes every run with matching formatting to be re-coalesced
# into one big run, so each boundary we need to keep is "pinned" by
# toggling a direct-character-formatting property on and back off, which
# forces a run split at that exact point without leaving any visible
# formatting difference behind.

$d = $word.ActiveDocument

# --- 1. Date field result -------------------------------------------------
$d.Content.Find.Execute("March 5, 2021", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "March 30, 2021", 2) | Out-Null

# --- 2. "the cutting-edge" -> "a cutting-edge", split into 3 runs --------
$findRange = $d.Content
$findRange.Find.Execute("the cutting-edge") | Out-Null
$theStart = $findRange.Start

# Replace "the" with "a" (3 chars -> 1 char).
$theRange = $d.Range($theStart, $theStart + 3)
$theRange.Text = "a"

# Pin "a" into its own run.
$aRange = $d.Range($theStart, $theStart + 1)
$aRange.Bold = 1
$aRange.Bold = 0

# Pin the existing "Coinbase" run so it doesn't get folded into the new
# trailing run.
$coinRange = $d.Content
$coinRange.Start = $theStart + 1
$coinRange.Find.Execute("Coinbase") | Out-Null
$coinRange.Bold = 1
$coinRange.Bold = 0

# Pin the ". " run that follows "Coinbase".
$periodRange = $d.Range($coinRange.End, $coinRange.End + 2)
$periodRange.Bold = 1
$periodRange.Bold = 0

# Pin the "As a data analyst..." run.
$asRange = $d.Content
$asRange.Start = $periodRange.End
$asRange.Find.Execute("As a data analyst at your company, I would be able " + `
    "to apply my experience to provide data insights that help businesses " + `
    "make informed decisions and reduce risks. ") | Out-Null
$asRange.Bold = 1
$asRange.Bold = 0

# Pin the closing "Thank you..." run.
$thankRange = $d.Content
$thankRange.Start = $asRange.End
$thankRange.Find.Execute("Thank you for taking the time to review my " + `
    "application; I look forward to an opportunity to learn more about " + `
    "the position and to further discuss my relevant skills and experience. ") | Out-Null
$thankRange.Bold = 1
$thankRange.Bold = 0
